$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, and the
# (now-unwanted) blank paragraph immediately preceding it.
$findRng = $d.Content
$findRng.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$verPara = $findRng.Paragraphs(1)
$emptyPara = $verPara.Previous()
$delStart = $emptyPara.Range.Start

# Locate the end of the "© 2020 . Contact: ..." paragraph that follows it.
$findRng2 = $d.Content
$findRng2.Find.Execute("Powered by Jekyll and Github pages", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$copyrightPara = $findRng2.Paragraphs(1)
$delEnd = $copyrightPara.Range.End

# Remove the blank paragraph + "Ver no Jupiter..." paragraph + "© 2020..."
# paragraph in one go, leaving the "LOB1037: ..." paragraph before them, and
# the blank/page-break paragraphs after them, untouched.
$r = $d.Range($delStart, $delEnd)
$r.Delete()
